# Fill in the new "carrier" column (D) for the practice rows (2-5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Fill in the new "pair_kind" column (J) for rows 6-9 (video pairs then audio pairs)
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Populate the previously-empty rows 14-21 with kind (C) and carrier (D) data
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
